$wb = $excel.ActiveWorkbook

# 1) Lead_Config20_CL: C2 500 -> 300
$ws = $wb.Worksheets.Item("Lead_Config20_CL")
$ws.Activate()
$ws.Range("C2").Value = 300
$ws.Range("C2").Select()

# 2) Lead_Config20_1_CL: C2 "Jaipur" -> "Mumbai"
$ws = $wb.Worksheets.Item("Lead_Config20_1_CL")
$ws.Activate()
$ws.Range("C2").Value = "Mumbai"
$ws.Range("A2").Select()

# 3) Lead_Config21_CL: C2 "India" -> "afghanistan"
$ws = $wb.Worksheets.Item("Lead_Config21_CL")
$ws.Activate()
$ws.Range("C2").Value = "afghanistan"
$ws.Range("C2").Select()
